# Add 2022-Q3 data
# -----------------------------------------------------------------
# 1. Insert a new worksheet (a copy of the existing "2022-Q2" sheet,
#    so it inherits identical column layout / styles) right after the
#    "总计" (total) sheet, then rename it to "2022-Q3" and fill in the
#    new quarter's fund-holding data.
# 2. Update the "总计" summary sheet: shift the existing quarterly
#    rows down by one and insert the new 2022-Q3 totals at the top.
# -----------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item(2)

# --- create the new "2022-Q3" sheet right after "总计" -------------
$q2Sheet.Copy($null, $total)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Row 2: 001959 / 华商乐享互联灵活配置混合A
$q3Sheet.Cells.Item(2,1).Value = 0
$q3Sheet.Cells.Item(2,2).Value = "'001959"
$q3Sheet.Cells.Item(2,3).Value = "华商乐享互联灵活配置混合A"
$q3Sheet.Cells.Item(2,4).Value = "'4.62"
$q3Sheet.Cells.Item(2,5).Value = "'93.28"
$q3Sheet.Cells.Item(2,6).Value = "'3.74"
$q3Sheet.Cells.Item(2,7).Value = "'0.1728"
$q3Sheet.Cells.Item(2,8).Value = 2

# Row 3: 013142 / 华商乐享互联灵活配置混合C
$q3Sheet.Cells.Item(3,1).Value = 1
$q3Sheet.Cells.Item(3,2).Value = "'013142"
$q3Sheet.Cells.Item(3,3).Value = "华商乐享互联灵活配置混合C"
$q3Sheet.Cells.Item(3,4).Value = "'1.08"
$q3Sheet.Cells.Item(3,5).Value = "'93.28"
$q3Sheet.Cells.Item(3,6).Value = "'3.74"
$q3Sheet.Cells.Item(3,7).Value = "'0.0404"
$q3Sheet.Cells.Item(3,8).Value = 2

# --- update the "总计" sheet ---------------------------------------
# Shift existing rows 2..8 down to rows 3..9 (bottom-up so we never
# clobber a row before it has been read).
$total.Cells.Item(9,2).Value = $total.Cells.Item(8,2).Value2
$total.Cells.Item(9,3).Value = $total.Cells.Item(8,3).Value2
$total.Cells.Item(9,4).Value = $total.Cells.Item(8,4).Value2

$total.Cells.Item(8,2).Value = $total.Cells.Item(7,2).Value2
$total.Cells.Item(8,3).Value = $total.Cells.Item(7,3).Value2
$total.Cells.Item(8,4).Value = $total.Cells.Item(7,4).Value2

$total.Cells.Item(7,2).Value = $total.Cells.Item(6,2).Value2
$total.Cells.Item(7,3).Value = $total.Cells.Item(6,3).Value2
$total.Cells.Item(7,4).Value = $total.Cells.Item(6,4).Value2

$total.Cells.Item(6,2).Value = $total.Cells.Item(5,2).Value2
$total.Cells.Item(6,3).Value = $total.Cells.Item(5,3).Value2
$total.Cells.Item(6,4).Value = $total.Cells.Item(5,4).Value2

$total.Cells.Item(5,2).Value = $total.Cells.Item(4,2).Value2
$total.Cells.Item(5,3).Value = $total.Cells.Item(4,3).Value2
$total.Cells.Item(5,4).Value = $total.Cells.Item(4,4).Value2

$total.Cells.Item(4,2).Value = $total.Cells.Item(3,2).Value2
$total.Cells.Item(4,3).Value = $total.Cells.Item(3,3).Value2
$total.Cells.Item(4,4).Value = $total.Cells.Item(3,4).Value2

$total.Cells.Item(3,2).Value = $total.Cells.Item(2,2).Value2
$total.Cells.Item(3,3).Value = $total.Cells.Item(2,3).Value2
$total.Cells.Item(3,4).Value = $total.Cells.Item(2,4).Value2

# New row 2: 2022-Q3 totals
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 2
$total.Cells.Item(2,4).Value = 0.21

# Re-number the index column (A) 0..7 top to bottom
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(6,1).Value = 4
$total.Cells.Item(7,1).Value = 5
$total.Cells.Item(8,1).Value = 6
$total.Cells.Item(9,1).Value = 7

# Row 9 is brand new -- copy the A-column formatting from row 8 so it
# matches the other index cells (bold border style).
$total.Cells.Item(8,1).Copy()
$total.Cells.Item(9,1).PasteSpecial(-4122)
